$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume-change (column E) values
# for the rows refreshed by this data sync. Column D text values are
# prefixed with an apostrophe so Excel stores them as text (matching the
# original inlineStr cell type) instead of auto-converting to numbers.
$ws.Range("D2").Value = "'24.990.82"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "'1.702.93"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'316.42"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "'0.3972"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").Value = "'0.4028"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'1.469"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").Value = "'1.005"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'0.08809"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'25.95"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "'7.461"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "'7.966"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "'0.00001350"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'1.711.02"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'96.36"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "'0.07206"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'20.62"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "'7.349"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'14.38"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'24.978.57"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "'2.350"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'2.934"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("D27").Value = "'23.74"
$ws.Range("E27").Value = "  +4.97%  "
$ws.Range("D28").Value = "'6.145"
$ws.Range("E28").Value = "  +14.52%  "
$ws.Range("D29").Value = "'162.15"
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("D30").Value = "'149.87"
$ws.Range("E30").Value = "  +8.26%  "
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Value = "'2.649"
$ws.Range("E32").Value = "  +25.19%  "
$ws.Range("D33").Value = "'1.899.42"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "'0.08555"
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").Value = "'0.03136"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("D36").Value = "'1.041"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "'7.134"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "'0.2854"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").Value = "'10.89"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "'0.09545"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").Value = "'0.8244"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("D42").Value = "'13.97"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "'1.479"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "'2.679"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").Value = "'0.7385"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "'4.257"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "'1.404"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "'0.08742"
$ws.Range("E49").Value = "  +8.31%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'139.21"
$ws.Range("E51").Value = "  -0.10%  "
